$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "30.675.15"
$ws.Range("E2").Value = "  -0.83%  "
Set-TextCell "D3" "1.889.66"
$ws.Range("E3").Value = "  -1.01%  "
Set-TextCell "D4" "1.000"
Set-TextCell "D5" "237.22"
$ws.Range("E5").Value = "  -3.61%  "
Set-TextCell "D6" "1.000"
$ws.Range("E6").Value = "  -0.02%  "
Set-TextCell "D7" "0.4879"
$ws.Range("E7").Value = "  -2.43%  "
Set-TextCell "D8" "0.2925"
$ws.Range("E8").Value = "  -2.09%  "
Set-TextCell "D9" "0.06681"
$ws.Range("E9").Value = "  -2.37%  "
Set-TextCell "D10" "1.890.77"
$ws.Range("E10").Value = "  -0.93%  "
Set-TextCell "D11" "16.72"
$ws.Range("E11").Value = "  -1.74%  "
Set-TextCell "D12" "0.07235"
$ws.Range("E12").Value = "  -1.36%  "
Set-TextCell "D13" "89.32"
$ws.Range("E13").Value = "  -2.45%  "
Set-TextCell "D14" "5.010"
$ws.Range("E14").Value = "  -1.96%  "
Set-TextCell "D15" "0.6650"
$ws.Range("E15").Value = "  -2.29%  "
Set-TextCell "D16" "30.627.54"
$ws.Range("E16").Value = "  -0.92%  "
Set-TextCell "D17" "0.000007917"
$ws.Range("E17").Value = "  -1.67%  "
Set-TextCell "D18" "1.000"
Set-TextCell "D19" "13.02"
$ws.Range("E19").Value = "  -2.05%  "
Set-TextCell "D20" "2.134.69"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("E21").Value = "  +0.28%  "
Set-TextCell "D22" "4.747"
$ws.Range("E22").Value = "  -2.81%  "
Set-TextCell "D23" "192.32"
$ws.Range("E23").Value = "  +5.02%  "
Set-TextCell "D24" "6.075"
$ws.Range("E24").Value = "  -0.51%  "
Set-TextCell "D25" "9.312"
$ws.Range("E25").Value = "  -0.61%  "
Set-TextCell "D26" "159.77"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("E27").Value = "  -2.14%  "
Set-TextCell "D28" "1.833"
$ws.Range("E28").Value = "  -6.01%  "
Set-TextCell "D29" "1.402"
$ws.Range("E29").Value = "  +0.54%  "
Set-TextCell "D30" "4.266"
$ws.Range("E30").Value = "  -2.02%  "
Set-TextCell "D31" "0.09031"
$ws.Range("E31").Value = "  +0.31%  "
Set-TextCell "D32" "3.941"
$ws.Range("E32").Value = "  -3.17%  "
Set-TextCell "D33" "0.05208"
$ws.Range("E33").Value = "  -1.31%  "
Set-TextCell "D34" "0.7330"
$ws.Range("E34").Value = "  -2.01%  "
Set-TextCell "D35" "1.087"
$ws.Range("E35").Value = "  -4.59%  "
Set-TextCell "D36" "2.684"
$ws.Range("E36").Value = "  +0.58%  "
Set-TextCell "D37" "0.01823"
$ws.Range("E37").Value = "  -6.21%  "
$ws.Range("E38").Value = "  -2.57%  "
Set-TextCell "D39" "0.9261"
$ws.Range("E39").Value = "  -1.25%  "
Set-TextCell "D40" "2.053"
$ws.Range("E40").Value = "  -6.17%  "
Set-TextCell "D41" "0.4408"
$ws.Range("E41").Value = "  -0.14%  "
Set-TextCell "D42" "104.58"
$ws.Range("E42").Value = "  -1.83%  "
Set-TextCell "D43" "0.9995"
$ws.Range("E43").Value = "  -0.12%  "
Set-TextCell "D44" "5.743"
$ws.Range("E44").Value = "  -1.97%  "
Set-TextCell "D45" "0.1338"
$ws.Range("E45").Value = "  -0.64%  "
Set-TextCell "D46" "7.355"
$ws.Range("E46").Value = "  -5.56%  "
Set-TextCell "D47" "0.4157"
$ws.Range("E47").Value = "  +5.88%  "
Set-TextCell "D48" "0.05831"
$ws.Range("E48").Value = "  -0.28%  "
Set-TextCell "D49" "8.734"
$ws.Range("E49").Value = "  +1.88%  "
Set-TextCell "D50" "1.408"
$ws.Range("E50").Value = "  +1.13%  "
Set-TextCell "D51" "33.25"
$ws.Range("E51").Value = "  -0.45%  "
